# Add five new "foto" columns (foto4..foto8) to the header row of the
# products sheet, mirroring the existing foto0..foto3 columns (I1:L1).
# New header cells M1:Q1 get the same value/border/font formatting as the
# existing header cells, and the active selection moves to the newly
# added column Q1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "foto4"
$ws.Range("N1").Value = "foto5"
$ws.Range("O1").Value = "foto6"
$ws.Range("P1").Value = "foto7"
$ws.Range("Q1").Value = "foto8"

# Copy the header formatting (font/border/alignment) from an existing
# header cell onto the new ones, same as Excel's "fill right" / paste
# formats behaviour would produce.
[void]$ws.Range("L1").Copy()
[void]$ws.Range("M1:Q1").PasteSpecial(-4122)

# Move the view / selection to the newly added column, like the author
# scrolled over to inspect the new headers after pasting them.
[void]$ws.Range("Q1").Select()
